$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = '550'
$ws.Range('C2').Value = 'No, FLASCO does not encompass community sites. FLASCO is primarily focused on education and advocacy for medical oncologists and other professionals in the field of oncology, rather than community-based practices.'
$ws.Range('D2').Value = 'No, FLASCO is an organization focused on education, research, and advocacy for clinical oncology professionals, rather than directly influencing policy.'
$ws.Range('E2').Value = 'Yes, FLASCO provides engagement opportunity with leadership. FLASCO offers networking events, committee involvement, and opportunities to interact with leaders in the field of clinical oncology.'
$ws.Range('F2').Value = 'No, FLASCO does not provide support for clinical trial recruitment. FLASCO''s primary focus is on education and advocacy for clinical oncologists in Florida, rather than directly assisting with recruitment for trials.'
$ws.Range('G2').Value = 'No, FLASCO does not provide engagement opportunity with payors. FLASCO''s focus is on supporting oncology professionals in Florida rather than facilitating interactions with payors.'
$ws.Range('H2').Value = 'Yes, FLASCO includes area experts on its board. The organization consists of Florida-based oncology professionals who are considered experts in the field.'
$ws.Range('I2').Value = 'Yes, FLASCO is involved in therapeutic research collaborations. FLASCO works with various organizations and institutions on clinical trials and research studies to advance cancer treatments.'
$ws.Range('J2').Value = 'Yes, FLASCO includes top therapeutic area experts on its board. This can be justified by the fact that FLASCO is a professional organization for medical oncologists in Florida, so it would make sense for the board to include experts in the field.'

$ws.Range('B3').Value = '400'
$ws.Range('C3').Value = 'No, GASCO does not encompass community sites, as it is specifically focused on clinical oncology.'
$ws.Range('D3').Value = 'No, GASCO does not have direct influence on state or local policy as it is primarily a professional organization focused on clinical oncology practices and advocacy.'
$ws.Range('E3').Value = 'Yes, GASCO provides engagement opportunities with leadership. GASCO hosts various events, meetings, and conferences where members can interact with leadership and get involved in decision-making processes.'
$ws.Range('F3').Value = 'Yes, GASCO provides support for clinical trial recruitment through education and outreach efforts to oncologists and patients.'
$ws.Range('G3').Value = 'No, GASCO does not provide engagement opportunities with payors. GASCO primarily focuses on providing education, networking, and advocacy for clinical oncologists in Georgia, rather than facilitating interactions with payors.'
$ws.Range('H3').Value = 'Yes, Gasco includes area experts on its board as they ensure the organization benefits from a wide range of knowledge and expertise in the field of clinical oncology.'
$ws.Range('I3').Value = 'Yes, GASCO is involved in therapeutic research collaborations. GASCO partners with various organizations to advance research in oncology.'
$ws.Range('J3').Value = 'No, GASCO does not have top therapeutic area experts on its board. The organization primarily consists of clinical oncologists and other professionals in oncology-related fields.'

$ws.Range('B4').Value = '550'
$ws.Range('C4').Value = 'No, IOS does not encompass community sites. The organization is dedicated to oncology practices and professionals in Indiana, typically focusing on academic or hospital-based settings.'
$ws.Range('D4').Value = 'No, IOS is primarily focused on providing education and resources for oncology professionals rather than directly influencing policy.'
$ws.Range('E4').Value = 'Yes, IOS provides leadership engagement opportunities for members. The organization offers networking events, educational programs, and mentorship opportunities that allow members to connect with leaders in the field of oncology.'
$ws.Range('F4').Value = 'No, IOS does not provide support for clinical trial recruitment. IOS is primarily focused on providing education, networking, and advocacy for oncology professionals in Indiana, rather than directly assisting in patient recruitment for clinical trials.'
$ws.Range('G4').Value = 'No, IOS does not provide engagement opportunities with payors. IOS is a professional organization focused on oncology education and support, rather than payor negotiations.'
$ws.Range('H4').Value = 'No,IOS does not include area experts on its board. IOS''s board is comprised of oncologists and healthcare professionals who specialize in the field of oncology.'
$ws.Range('I4').Value = 'No, IOS is primarily focused on education and advocacy for oncologists in Indiana.'
$ws.Range('J4').Value = 'No,  The Indiana Oncology Society (IOS) does not include top therapeutic area experts on its board. , There is no public information confirming the presence of top experts on the board.'

$ws.Range('B5').Value = '100'
$ws.Range('C5').Value = 'No, because it is specific to Oncology and likely focuses on larger medical facilities.'
$ws.Range('D5').Value = 'No, limited membership and focus on oncology issues.'
$ws.Range('E5').Value = @'
Yes, 
The IOWA Oncology Society provides various opportunities for its members to engage with leadership through meetings, events, and committees.
'@
$ws.Range('F5').Value = 'Yes, they provide support for clinical trial recruitment as part of their mission to advance cancer treatment options for patients.'
$ws.Range('G5').Value = 'No, The IOWA Oncology Society is a professional organization focused on education and support for oncologists, not on engaging with payors.'
$ws.Range('H5').Value = 'No, The IOWA Oncology Society does not have area experts on its board. The society is composed of oncology professionals from various backgrounds.'
$ws.Range('I5').Value = 'Yes, IOWA Oncology Society is likely involved in therapeutic research collaborations as they are a professional organization focused on oncology, which often involves research collaborations.'
$ws.Range('J5').Value = 'Yes, the IOWA Oncology Society includes top therapeutic area experts on its board. The board members are actively involved in the field of oncology and bring vast experience and knowledge to the society''s activities and initiatives.'

$ws.Range('B6').Value = '340'
$ws.Range('C6').Value = 'No, MOASC does not encompass community sites, as it is an association focused on medical oncology in Southern California.'
$ws.Range('D6').Value = 'No, MOASC does not have a direct influence on state or local policy as it is primarily a professional organization focused on medical oncology practices and education.'
$ws.Range('E6').Value = 'No, MOASC does not provide engagement opportunity with leadership, as it is primarily focused on networking and education for medical oncology professionals.'
$ws.Range('F6').Value = 'No, MOASC does not provide support for clinical trial recruitment. , MOASC focuses primarily on education, advocacy, and networking for medical oncologists in Southern California.'
$ws.Range('G6').Value = 'No, MOASC primarily focuses on providing support, education, and advocacy for medical oncologists in Southern California, rather than engaging with payors.'
$ws.Range('H6').Value = 'Yes, the MOASC board includes area experts as they typically consist of leading oncologists and healthcare professionals in Southern California.'
$ws.Range('I6').Value = 'Yes, MOASC is involved in therapeutic research collaborations. MOASC works closely with pharmaceutical companies and academic institutions to conduct clinical research studies aimed at advancing cancer treatment options for patients.'
$ws.Range('J6').Value = @'
Yes, 
The MOASC board includes top therapeutic area experts who are leaders in the field of medical oncology in Southern California.
'@

$ws.Range('K4').Value = 'Midwest'
$ws.Range('K5').Value = 'Midwest'
